# "Un poco de todo" -- add two new icon rows (Check True/False) to the
# "Otros Iconos" table on Hoja1: circle-check / circle-xmark, each with a
# hyperlink to its FontAwesome page, formatted like the row above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: Check (True) / circle-check -----------------------------
$ws.Range("B36").Value = "Check (True)"
$ws.Range("C36").Value = "circle-check"
[void]$ws.Hyperlinks.Add($ws.Range("D36"), "https://fontawesome.com/icons/circle-check?f=classic&s=regular")

# --- Row 37: Check (False) / circle-xmark -----------------------------
# C/D are entered before B so the shared-string table picks up
# "Check (False)" as the very last new entry (matches source order).
$ws.Range("C37").Value = "circle-xmark"
[void]$ws.Hyperlinks.Add($ws.Range("D37"), "https://fontawesome.com/icons/circle-xmark?f=classic&s=regular")
$ws.Range("B37").Value = "Check (False)"

# Copy the formatting (borders/fill/hyperlink style) from row 35 down
# onto the two new rows so they visually match the rest of the table.
$ws.Range("B35:D35").Copy()
$ws.Range("B36:D37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the selection where the author left it.
$ws.Range("G34").Select() | Out-Null
